# Datos de colaboradores DB y correccion en controlador
#
# Update the "Rol" (D) and "Disponibilidad" (F) values for the existing
# collaborator rows, and move the current selection to E6 (matching the
# saved cursor position in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 2 - Adrian Ramirez
$ws.Range("D2").Value = 12

# Row 3 - Celina Madrigal
$ws.Range("D3").Value = 13
$ws.Range("F3").Value = 2

# Row 4 - Maria Porras
$ws.Range("D4").Value = 14
$ws.Range("F4").Value = 3

# Restore the cursor / selection position recorded on last save
$ws.Activate()
$ws.Range("E6").Select() | Out-Null
